$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 575
$ws.Cells.Item(12, 9).Value = 150
$ws.Cells.Item(12, 11).Value = 150
$ws.Cells.Item(12, 13).Value = 20
$ws.Cells.Item(33, 8).Value = 508.56
$ws.Cells.Item(33, 9).Value = 361.69565
$ws.Cells.Item(33, 10).Value = 2197.5
$ws.Cells.Item(33, 11).Value = 361.69565
$ws.Cells.Item(33, 12).Value = 2197.5
$ws.Cells.Item(33, 13).Value = -132.69565
$ws.Cells.Item(33, 14).Value = -2655.5
$ws.Cells.Item(34, 8).Value = 1982.6666
$ws.Cells.Item(34, 9).Value = 1856.4286
$ws.Cells.Item(34, 10).Value = 2424.5
$ws.Cells.Item(34, 11).Value = 1856.4286
$ws.Cells.Item(34, 12).Value = 2424.5
$ws.Cells.Item(34, 13).Value = -1653.4286
$ws.Cells.Item(34, 14).Value = -2830.5
$ws.Cells.Item(36, 8).Value = 1982.6666
$ws.Cells.Item(36, 9).Value = 1856.4286
$ws.Cells.Item(36, 10).Value = 2424.5
$ws.Cells.Item(36, 11).Value = 1856.4286
$ws.Cells.Item(36, 12).Value = 2424.5
$ws.Cells.Item(36, 13).Value = -1141.4286
$ws.Cells.Item(36, 14).Value = -3854.5
$ws.Cells.Item(53, 8).Value = 443.36365
$ws.Cells.Item(53, 9).Value = 183.2
$ws.Cells.Item(53, 11).Value = 183.2
$ws.Cells.Item(53, 13).Value = 453.8
$ws.Cells.Item(62, 8).Value = 8409690
$ws.Cells.Item(62, 9).Value = 5048.875
$ws.Cells.Item(62, 11).Value = 5048.875
$ws.Cells.Item(62, 13).Value = -4424.875
$ws.Cells.Item(65, 8).Value = 8409690
$ws.Cells.Item(65, 9).Value = 5048.875
$ws.Cells.Item(65, 11).Value = 25244.375
$ws.Cells.Item(65, 13).Value = -22124.375
$ws.Cells.Item(92, 8).Value = 1203.3125
$ws.Cells.Item(92, 9).Value = 1111.3846
$ws.Cells.Item(92, 10).Value = 1601.6666
$ws.Cells.Item(92, 11).Value = 1111.3846
$ws.Cells.Item(92, 12).Value = 1601.6666
$ws.Cells.Item(92, 13).Value = 136.6153999999999
$ws.Cells.Item(92, 14).Value = -4097.6666
$ws.Cells.Item(96, 8).Value = 1653.3667
$ws.Cells.Item(96, 9).Value = 1669.0555
$ws.Cells.Item(96, 10).Value = 1629.8334
$ws.Cells.Item(96, 11).Value = 5007.166499999999
$ws.Cells.Item(96, 12).Value = 4889.5002
$ws.Cells.Item(96, 13).Value = -3634.166499999999
$ws.Cells.Item(96, 14).Value = -7635.5002
$ws.Cells.Item(98, 8).Value = 4889.636
$ws.Cells.Item(98, 9).Value = 4889.636
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 4889.636
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = -3391.636
$ws.Cells.Item(98, 14).ClearContents()
$ws.Cells.Item(99, 8).Value = 736.8461
$ws.Cells.Item(99, 9).Value = 223.75
$ws.Cells.Item(99, 11).Value = 671.25
$ws.Cells.Item(99, 13).Value = 826.75
$ws.Cells.Item(107, 8).Value = 2019.8572
$ws.Cells.Item(107, 9).Value = 690.8333
$ws.Cells.Item(107, 11).Value = 690.8333
$ws.Cells.Item(107, 13).Value = 1229.1667
$ws.Cells.Item(116, 8).Value = 6114.8335
$ws.Cells.Item(116, 9).Value = 4900
$ws.Cells.Item(116, 11).Value = 4900
$ws.Cells.Item(116, 13).Value = -1458
$ws.Cells.Item(122, 8).Value = 4889.636
$ws.Cells.Item(122, 9).Value = 4889.636
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 14668.908
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -12218.908
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 1905.6809
$ws.Cells.Item(132, 9).Value = 1729.525
$ws.Cells.Item(132, 11).Value = 5188.575000000001
$ws.Cells.Item(132, 13).Value = -2658.575000000001
$ws.Cells.Item(141, 8).Value = 11331.786
$ws.Cells.Item(141, 9).Value = 11895.77
$ws.Cells.Item(141, 11).Value = 35687.31
$ws.Cells.Item(141, 13).Value = -30507.31

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).ClearContents()
$ws.Cells.Item(2, 14).ClearContents()
$ws.Cells.Item(18, 8).Value = 1470
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).ClearContents()
$ws.Cells.Item(32, 8).Value = 11364870
$ws.Cells.Item(32, 9).Value = 12346552
$ws.Cells.Item(32, 10).Value = 5398.7144
$ws.Cells.Item(32, 11).Value = 12346552
$ws.Cells.Item(32, 12).Value = 5398.7144
$ws.Cells.Item(32, 13).Value = -12346265
$ws.Cells.Item(32, 14).Value = -5972.7144
$ws.Cells.Item(45, 8).Value = 2523.1538
$ws.Cells.Item(45, 9).Value = 1949.25
$ws.Cells.Item(45, 10).Value = 2778.2222
$ws.Cells.Item(45, 11).Value = 1949.25
$ws.Cells.Item(45, 12).Value = 2778.2222
$ws.Cells.Item(45, 13).Value = -1572.25
$ws.Cells.Item(45, 14).Value = -3532.2222
$ws.Cells.Item(61, 8).Value = 23816876
$ws.Cells.Item(61, 9).Value = 23816876
$ws.Cells.Item(61, 11).Value = 23816876
$ws.Cells.Item(61, 13).Value = -23816664
$ws.Cells.Item(63, 8).Value = 6469.5
$ws.Cells.Item(63, 9).Value = 1563.3334
$ws.Cells.Item(63, 11).Value = 1563.3334
$ws.Cells.Item(63, 13).Value = -877.3334
$ws.Cells.Item(66, 8).Value = 6469.5
$ws.Cells.Item(66, 9).Value = 1563.3334
$ws.Cells.Item(66, 11).Value = 7816.666999999999
$ws.Cells.Item(66, 13).Value = -4384.666999999999
$ws.Cells.Item(107, 8).Value = 51666.332
$ws.Cells.Item(107, 9).Value = 73000
$ws.Cells.Item(107, 10).Value = 40999.5
$ws.Cells.Item(107, 11).Value = 73000
$ws.Cells.Item(107, 12).Value = 40999.5
$ws.Cells.Item(107, 13).Value = -69160
$ws.Cells.Item(107, 14).Value = -48679.5
$ws.Cells.Item(109, 8).Value = 60377
$ws.Cells.Item(109, 10).Value = 60377
$ws.Cells.Item(109, 12).Value = 60377
$ws.Cells.Item(109, 14).Value = -63151
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).ClearContents()
$ws.Cells.Item(116, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 29424542
$ws.Cells.Item(132, 9).Value = 5315.8213
$ws.Cells.Item(132, 10).Value = 166714270
$ws.Cells.Item(132, 11).Value = 15947.4639
$ws.Cells.Item(132, 12).Value = 500142810
$ws.Cells.Item(132, 13).Value = -13417.4639
$ws.Cells.Item(132, 14).Value = -500147870
$ws.Cells.Item(135, 8).Value = 79998
$ws.Cells.Item(135, 10).Value = 79998
$ws.Cells.Item(135, 12).Value = 79998
$ws.Cells.Item(135, 14).Value = -90138
$ws.Cells.Item(136, 8).Value = 23816876
$ws.Cells.Item(136, 9).Value = 23816876
$ws.Cells.Item(136, 11).Value = 71450628
$ws.Cells.Item(136, 13).Value = -71448078

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(3, 14).ClearContents()
$ws.Cells.Item(24, 8).Value = 8997.5
$ws.Cells.Item(24, 9).Value = 8000
$ws.Cells.Item(24, 10).Value = 9995
$ws.Cells.Item(24, 11).Value = 8000
$ws.Cells.Item(24, 12).Value = 9995
$ws.Cells.Item(24, 13).Value = -7765
$ws.Cells.Item(24, 14).Value = -10465
$ws.Cells.Item(30, 8).Value = 9999
$ws.Cells.Item(30, 10).Value = 9999
$ws.Cells.Item(30, 12).Value = 9999
$ws.Cells.Item(30, 14).Value = -10249
$ws.Cells.Item(86, 8).Value = 17170.125
$ws.Cells.Item(86, 9).Value = 18339.428
$ws.Cells.Item(86, 10).Value = 8985
$ws.Cells.Item(86, 11).Value = 18339.428
$ws.Cells.Item(86, 12).Value = 8985
$ws.Cells.Item(86, 13).Value = -17216.428
$ws.Cells.Item(86, 14).Value = -11231
$ws.Cells.Item(89, 8).Value = 17170.125
$ws.Cells.Item(89, 9).Value = 18339.428
$ws.Cells.Item(89, 10).Value = 8985
$ws.Cells.Item(89, 11).Value = 91697.14
$ws.Cells.Item(89, 12).Value = 44925
$ws.Cells.Item(89, 13).Value = -86081.14
$ws.Cells.Item(89, 14).Value = -56157
$ws.Cells.Item(134, 8).Value = 3377.4666
$ws.Cells.Item(134, 9).Value = 2333.0908
$ws.Cells.Item(134, 10).Value = 6249.5
$ws.Cells.Item(134, 11).Value = 6999.2724
$ws.Cells.Item(134, 12).Value = 18748.5
$ws.Cells.Item(134, 13).Value = -4464.2724
$ws.Cells.Item(134, 14).Value = -23818.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 56821620
$ws.Cells.Item(31, 9).Value = 2648.818
$ws.Cells.Item(31, 10).Value = 113640590
$ws.Cells.Item(31, 11).Value = 2648.818
$ws.Cells.Item(31, 12).Value = 113640590
$ws.Cells.Item(31, 13).Value = -2353.818
$ws.Cells.Item(31, 14).Value = -113641180
$ws.Cells.Item(34, 8).Value = 56821620
$ws.Cells.Item(34, 9).Value = 2648.818
$ws.Cells.Item(34, 10).Value = 113640590
$ws.Cells.Item(34, 11).Value = 2648.818
$ws.Cells.Item(34, 12).Value = 113640590
$ws.Cells.Item(34, 13).Value = -2446.818
$ws.Cells.Item(34, 14).Value = -113640994
$ws.Cells.Item(107, 8).Value = 980.5
$ws.Cells.Item(107, 9).Value = 980.5
$ws.Cells.Item(107, 11).Value = 980.5
$ws.Cells.Item(107, 13).Value = 939.5
$ws.Cells.Item(132, 8).Value = 1808.4193
$ws.Cells.Item(132, 9).Value = 1768.2333
$ws.Cells.Item(132, 11).Value = 5304.699900000001
$ws.Cells.Item(132, 13).Value = -2774.699900000001
$ws.Cells.Item(134, 8).Value = 3227.182
$ws.Cells.Item(134, 9).Value = 2581.2354
$ws.Cells.Item(134, 10).Value = 5423.4
$ws.Cells.Item(134, 11).Value = 7743.706200000001
$ws.Cells.Item(134, 12).Value = 16270.2
$ws.Cells.Item(134, 13).Value = -5208.706200000001
$ws.Cells.Item(134, 14).Value = -21340.2
$ws.Cells.Item(141, 8).Value = 196775.33
$ws.Cells.Item(141, 10).Value = 280163
$ws.Cells.Item(141, 12).Value = 280163
$ws.Cells.Item(141, 14).Value = -290523

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 12345726
$ws.Cells.Item(7, 9).Value = 15432135
$ws.Cells.Item(7, 10).Value = 90
$ws.Cells.Item(7, 11).Value = 46296405
$ws.Cells.Item(7, 12).Value = 270
$ws.Cells.Item(7, 13).Value = -46296293
$ws.Cells.Item(7, 14).Value = -494
$ws.Cells.Item(10, 8).Value = 1460.9375
$ws.Cells.Item(10, 9).Value = 155.85715
$ws.Cells.Item(10, 10).Value = 2476
$ws.Cells.Item(10, 11).Value = 467.57145
$ws.Cells.Item(10, 12).Value = 7428
$ws.Cells.Item(10, 13).Value = -328.57145
$ws.Cells.Item(10, 14).Value = -7706
$ws.Cells.Item(29, 8).Value = 6619.375
$ws.Cells.Item(29, 9).Value = 389
$ws.Cells.Item(29, 10).Value = 12849.75
$ws.Cells.Item(29, 11).Value = 1167
$ws.Cells.Item(29, 12).Value = 38549.25
$ws.Cells.Item(29, 13).Value = -890
$ws.Cells.Item(29, 14).Value = -39103.25
$ws.Cells.Item(86, 8).Value = 1399.7
$ws.Cells.Item(86, 10).Value = 1997.5
$ws.Cells.Item(86, 12).Value = 5992.5
$ws.Cells.Item(86, 14).Value = -8364.5
$ws.Cells.Item(89, 8).Value = 1399.7
$ws.Cells.Item(89, 10).Value = 1997.5
$ws.Cells.Item(89, 12).Value = 17977.5
$ws.Cells.Item(89, 14).Value = -29833.5
$ws.Cells.Item(95, 8).Value = 11000
$ws.Cells.Item(95, 10).Value = 11000
$ws.Cells.Item(95, 12).Value = 33000
$ws.Cells.Item(95, 14).Value = -37118
$ws.Cells.Item(98, 8).Value = 508.125
$ws.Cells.Item(98, 9).Value = 493
$ws.Cells.Item(98, 10).Value = 517.2
$ws.Cells.Item(98, 11).Value = 1479
$ws.Cells.Item(98, 12).Value = 1551.6
$ws.Cells.Item(98, 13).Value = 19
$ws.Cells.Item(98, 14).Value = -4547.6
$ws.Cells.Item(120, 8).Value = 20902.727
$ws.Cells.Item(120, 9).Value = 13732.5
$ws.Cells.Item(120, 11).Value = 41197.5
$ws.Cells.Item(120, 13).Value = -36359.5
$ws.Cells.Item(121, 8).Value = 949.5
$ws.Cells.Item(121, 9).Value = 949
$ws.Cells.Item(121, 10).Value = 950
$ws.Cells.Item(121, 11).Value = 2847
$ws.Cells.Item(121, 12).Value = 2850
$ws.Cells.Item(121, 13).Value = -1537
$ws.Cells.Item(121, 14).Value = -5470
$ws.Cells.Item(138, 8).Value = 8013
$ws.Cells.Item(138, 9).Value = 8013
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 24039
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 13).Value = -18899
$ws.Cells.Item(138, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 1160.6
$ws.Cells.Item(3, 10).Value = 1334.3334
$ws.Cells.Item(3, 12).Value = 1334.3334
$ws.Cells.Item(3, 14).Value = -1566.3334
$ws.Cells.Item(70, 8).Value = 11814.5
$ws.Cells.Item(70, 9).Value = 11251.75
$ws.Cells.Item(70, 10).Value = 12377.25
$ws.Cells.Item(70, 11).Value = 11251.75
$ws.Cells.Item(70, 12).Value = 12377.25
$ws.Cells.Item(70, 13).Value = -10981.75
$ws.Cells.Item(70, 14).Value = -12917.25
$ws.Cells.Item(73, 8).Value = 11814.5
$ws.Cells.Item(73, 9).Value = 11251.75
$ws.Cells.Item(73, 10).Value = 12377.25
$ws.Cells.Item(73, 11).Value = 11251.75
$ws.Cells.Item(73, 12).Value = 12377.25
$ws.Cells.Item(73, 13).Value = -10315.75
$ws.Cells.Item(73, 14).Value = -14249.25
$ws.Cells.Item(80, 8).Value = 9907.416999999999
$ws.Cells.Item(80, 9).Value = 8710
$ws.Cells.Item(80, 10).Value = 13499.667
$ws.Cells.Item(80, 11).Value = 8710
$ws.Cells.Item(80, 12).Value = 13499.667
$ws.Cells.Item(80, 13).Value = -7712
$ws.Cells.Item(80, 14).Value = -15495.667
$ws.Cells.Item(83, 8).Value = 9907.416999999999
$ws.Cells.Item(83, 9).Value = 8710
$ws.Cells.Item(83, 10).Value = 13499.667
$ws.Cells.Item(83, 11).Value = 43550
$ws.Cells.Item(83, 12).Value = 67498.33499999999
$ws.Cells.Item(83, 13).Value = -38558
$ws.Cells.Item(83, 14).Value = -77482.33499999999
$ws.Cells.Item(93, 8).Value = 35800
$ws.Cells.Item(93, 9).Value = 26333.334
$ws.Cells.Item(93, 11).Value = 26333.334
$ws.Cells.Item(93, 13).Value = -24461.334
$ws.Cells.Item(97, 8).Value = 690.1539
$ws.Cells.Item(97, 9).Value = 726.05884
$ws.Cells.Item(97, 11).Value = 726.05884
$ws.Cells.Item(97, 13).Value = -230.05884
$ws.Cells.Item(98, 8).Value = 685
$ws.Cells.Item(98, 10).Value = 685
$ws.Cells.Item(98, 12).Value = 685
$ws.Cells.Item(98, 14).Value = -6675
$ws.Cells.Item(102, 8).Value = 1888.85
$ws.Cells.Item(102, 9).Value = 1421.6666
$ws.Cells.Item(102, 10).Value = 2859.1538
$ws.Cells.Item(102, 11).Value = 1421.6666
$ws.Cells.Item(102, 12).Value = 2859.1538
$ws.Cells.Item(102, 13).Value = 200.3334
$ws.Cells.Item(102, 14).Value = -6103.1538
$ws.Cells.Item(113, 8).Value = 3980.7273
$ws.Cells.Item(113, 9).Value = 2863
$ws.Cells.Item(113, 10).Value = 4399.875
$ws.Cells.Item(113, 11).Value = 2863
$ws.Cells.Item(113, 12).Value = 4399.875
$ws.Cells.Item(113, 13).Value = -693
$ws.Cells.Item(113, 14).Value = -8739.875
$ws.Cells.Item(122, 8).Value = 21741208
$ws.Cells.Item(122, 9).Value = 1943.125
$ws.Cells.Item(122, 10).Value = 71430950
$ws.Cells.Item(122, 11).Value = 5829.375
$ws.Cells.Item(122, 12).Value = 214292850
$ws.Cells.Item(122, 13).Value = -3379.375
$ws.Cells.Item(122, 14).Value = -214297750
$ws.Cells.Item(126, 8).Value = 6669443
$ws.Cells.Item(126, 9).Value = 2336.8948
$ws.Cells.Item(126, 10).Value = 18185354
$ws.Cells.Item(126, 11).Value = 7010.6844
$ws.Cells.Item(126, 12).Value = 54556062
$ws.Cells.Item(126, 13).Value = -4540.6844
$ws.Cells.Item(126, 14).Value = -54561002
$ws.Cells.Item(132, 8).Value = 12007.921
$ws.Cells.Item(132, 9).Value = 10069.758
$ws.Cells.Item(132, 10).Value = 24799.8
$ws.Cells.Item(132, 11).Value = 30209.274
$ws.Cells.Item(132, 12).Value = 74399.39999999999
$ws.Cells.Item(132, 13).Value = -27679.274
$ws.Cells.Item(132, 14).Value = -79459.39999999999
$ws.Cells.Item(134, 8).Value = 47492.1
$ws.Cells.Item(134, 10).Value = 47492.1
$ws.Cells.Item(134, 12).Value = 142476.3
$ws.Cells.Item(134, 14).Value = -147546.3
$ws.Cells.Item(136, 8).Value = 86099
$ws.Cells.Item(136, 10).Value = 86099
$ws.Cells.Item(136, 12).Value = 258297
$ws.Cells.Item(136, 14).Value = -263397

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 1193
$ws.Cells.Item(19, 9).Value = 50
$ws.Cells.Item(19, 10).Value = 1478.75
$ws.Cells.Item(19, 11).Value = 50
$ws.Cells.Item(19, 12).Value = 1478.75
$ws.Cells.Item(19, 13).Value = 120
$ws.Cells.Item(19, 14).Value = -1818.75
$ws.Cells.Item(35, 8).Value = 1667.6666
$ws.Cells.Item(35, 9).Value = 1667.6666
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 1667.6666
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -1331.6666
$ws.Cells.Item(35, 14).ClearContents()
$ws.Cells.Item(40, 8).Value = 6044.9355
$ws.Cells.Item(40, 9).Value = 5681.522
$ws.Cells.Item(40, 10).Value = 7089.75
$ws.Cells.Item(40, 11).Value = 5681.522
$ws.Cells.Item(40, 12).Value = 7089.75
$ws.Cells.Item(40, 13).Value = -5545.522
$ws.Cells.Item(40, 14).Value = -7361.75
$ws.Cells.Item(55, 8).Value = 390.18182
$ws.Cells.Item(55, 9).Value = 349.625
$ws.Cells.Item(55, 10).Value = 498.33334
$ws.Cells.Item(55, 11).Value = 349.625
$ws.Cells.Item(55, 12).Value = 498.33334
$ws.Cells.Item(55, 13).Value = -176.625
$ws.Cells.Item(55, 14).Value = -844.33334
$ws.Cells.Item(61, 8).Value = 3609.25
$ws.Cells.Item(61, 9).Value = 2839.1428
$ws.Cells.Item(61, 11).Value = 2839.1428
$ws.Cells.Item(61, 13).Value = -2637.1428
$ws.Cells.Item(99, 8).Value = 41892
$ws.Cells.Item(99, 10).Value = 47676
$ws.Cells.Item(99, 12).Value = 47676
$ws.Cells.Item(99, 14).Value = -53666
$ws.Cells.Item(100, 8).Value = 2284.4333
$ws.Cells.Item(100, 9).Value = 1982.9259
$ws.Cells.Item(100, 11).Value = 1982.9259
$ws.Cells.Item(100, 13).Value = -1441.9259
$ws.Cells.Item(102, 8).Value = 66365
$ws.Cells.Item(102, 10).Value = 62780.5
$ws.Cells.Item(102, 12).Value = 62780.5
$ws.Cells.Item(102, 14).Value = -69270.5
$ws.Cells.Item(107, 8).Value = 25000
$ws.Cells.Item(107, 9).Value = 25000
$ws.Cells.Item(107, 11).Value = 25000
$ws.Cells.Item(107, 13).Value = -23080
$ws.Cells.Item(113, 8).Value = 3609.25
$ws.Cells.Item(113, 9).Value = 2839.1428
$ws.Cells.Item(113, 11).Value = 2839.1428
$ws.Cells.Item(113, 13).Value = -669.1428000000001
$ws.Cells.Item(136, 8).Value = 1253267.2
$ws.Cells.Item(136, 9).Value = 2502124.8
$ws.Cells.Item(136, 10).Value = 4409.875
$ws.Cells.Item(136, 11).Value = 7506374.399999999
$ws.Cells.Item(136, 12).Value = 13229.625
$ws.Cells.Item(136, 13).Value = -7503824.399999999
$ws.Cells.Item(136, 14).Value = -18329.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 42396
$ws.Cells.Item(64, 10).Value = 44993.332
$ws.Cells.Item(64, 12).Value = 44993.332
$ws.Cells.Item(64, 14).Value = -45489.332
$ws.Cells.Item(67, 8).Value = 42396
$ws.Cells.Item(67, 10).Value = 44993.332
$ws.Cells.Item(67, 12).Value = 44993.332
$ws.Cells.Item(67, 14).Value = -46709.332
$ws.Cells.Item(81, 8).Value = 3661.64
$ws.Cells.Item(81, 9).Value = 3110.4783
$ws.Cells.Item(81, 10).Value = 10000
$ws.Cells.Item(81, 11).Value = 6220.9566
$ws.Cells.Item(81, 12).Value = 20000
$ws.Cells.Item(81, 13).Value = -5159.9566
$ws.Cells.Item(81, 14).Value = -22122
$ws.Cells.Item(84, 8).Value = 3661.64
$ws.Cells.Item(84, 9).Value = 3110.4783
$ws.Cells.Item(84, 10).Value = 10000
$ws.Cells.Item(84, 11).Value = 31104.783
$ws.Cells.Item(84, 12).Value = 100000
$ws.Cells.Item(84, 13).Value = -25800.783
$ws.Cells.Item(84, 14).Value = -110608
$ws.Cells.Item(103, 8).Value = 54916.668
$ws.Cells.Item(103, 10).Value = 54916.668
$ws.Cells.Item(103, 12).Value = 54916.668
$ws.Cells.Item(103, 14).Value = -57260.668
$ws.Cells.Item(106, 8).Value = 50251
$ws.Cells.Item(106, 9).Value = 50000
$ws.Cells.Item(106, 11).Value = 50000
$ws.Cells.Item(106, 13).Value = -48738
$ws.Cells.Item(107, 8).Value = 2163.158
$ws.Cells.Item(107, 9).Value = 1700
$ws.Cells.Item(107, 10).Value = 2800
$ws.Cells.Item(107, 11).Value = 5100
$ws.Cells.Item(107, 12).Value = 8400
$ws.Cells.Item(107, 13).Value = -3180
$ws.Cells.Item(107, 14).Value = -12240
$ws.Cells.Item(113, 8).Value = 1932.75
$ws.Cells.Item(113, 9).Value = 1932.75
$ws.Cells.Item(113, 11).Value = 5798.25
$ws.Cells.Item(113, 13).Value = -3628.25
$ws.Cells.Item(122, 8).Value = 6253759
$ws.Cells.Item(122, 9).Value = 3706.087
$ws.Cells.Item(122, 10).Value = 22226116
$ws.Cells.Item(122, 11).Value = 11118.261
$ws.Cells.Item(122, 12).Value = 66678348
$ws.Cells.Item(122, 13).Value = -8668.261
$ws.Cells.Item(122, 14).Value = -66683248
$ws.Cells.Item(126, 8).Value = 2693.0356
$ws.Cells.Item(126, 9).Value = 1953.6842
$ws.Cells.Item(126, 10).Value = 4253.8887
$ws.Cells.Item(126, 11).Value = 5861.0526
$ws.Cells.Item(126, 12).Value = 12761.6661
$ws.Cells.Item(126, 13).Value = -3391.0526
$ws.Cells.Item(126, 14).Value = -17701.6661
$ws.Cells.Item(132, 8).Value = 3063.1738
$ws.Cells.Item(132, 9).Value = 2747.8
$ws.Cells.Item(132, 11).Value = 8243.400000000001
$ws.Cells.Item(132, 13).Value = -5713.400000000001
$ws.Cells.Item(136, 8).Value = 1817.4166
$ws.Cells.Item(136, 9).Value = 1201.4
$ws.Cells.Item(136, 11).Value = 3604.2
$ws.Cells.Item(136, 13).Value = -1054.2
